$wb = $excel.ActiveWorkbook

# Add the new "Numeric To Categorical" worksheet after the last existing sheet
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add($null, $lastSheet)
$newSheet.Name = "Numeric To Categorical"

# Header row
$newSheet.Range("A1").Value = "Action"
$newSheet.Range("B1").Value = "Time"
$newSheet.Range("C1").Value = "Content"

# Data rows
$newSheet.Range("A2").Value = "Upload CSV"
$newSheet.Range("B2").Value = "5 min"
$newSheet.Range("C2").Value = "df = pd.read_csv('file.csv')"

$newSheet.Range("A3").Value = "Identify Numerics"
$newSheet.Range("B3").Value = "2 min"
$newSheet.Range("C3").Value = "df.select_dtypes(include=['int', 'float'])"

$newSheet.Range("A4").Value = "Convert to Category"
$newSheet.Range("B4").Value = "2 min"
$newSheet.Range("C4").Value = "df['column'] = df['column'].astype('category')"

$newSheet.Range("A5").Value = "Verify Changes"
$newSheet.Range("B5").Value = "1 min"
$newSheet.Range("C5").Value = "df.info() to check new dtypes"

$newSheet.Range("A6").Value = "Overall"
$newSheet.Range("B6").Value = "10 min"

# Formatting: bold header + totals row, regular weight for body rows, all 13pt
# (set font size before bold to avoid creating a transient "bold 12pt" style)
$newSheet.Range("A1:C5").Font.Size = 13
$newSheet.Range("A6:B6").Font.Size = 13
$newSheet.Range("A1:C1").Font.Bold = $true
$newSheet.Range("A6:B6").Font.Bold = $true

# Row heights to match the other sheets in the workbook
$newSheet.Range("A1:A6").RowHeight = 17

# Select A1 as the active cell on the new sheet
$newSheet.Range("A1").Select() | Out-Null
